# Add files via upload
#
# The author's edit (as captured by the OOXML diff):
#   1. Renamed the "target" value "proton" to "p" for every data row.
#   2. Selected the header row A1:K1 and made it bold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Find the "target" column from the header row instead of hard-coding it,
# then swap every "proton" entry in that column for "p".
$targetCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Value2 -eq "target") {
        $targetCol = $c
    }
}

if ($targetCol -gt 0) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $targetCol)
        if ($cell.Value2 -eq "proton") {
            $cell.Value2 = "p"
        }
    }
}

# Bold the header row and leave it selected, matching the author's edit.
$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $lastCol))
$headerRange.Font.Bold = $true
$headerRange.Select()
